# ---------------------------------------------------------------------------
# CryCompanywiseStockReport_1 - stock count correction pass
#
# This script re-applies a corrected stock count/value extract on top of the
# previously-loaded workbook:
#   1) Three pairs of adjacent item rows had been entered in the wrong order
#      (duplicate-looking SKUs with differing case/spacing in the product name
#      were transposed) - their Item Code / Name / Rate / MRP / Qty / Value are
#      swapped back into the correct rows.
#   2) A batch of rows had their issued/closing Quantity (col F) corrected
#      (mostly down by a small count); the row Value (col G, = Rate * Qty) is
#      recalculated to match.
#   3) Every "Sub Total:" row (col B) affected by the above is recalculated as
#      the sum of the Value column for its item group, and the workbook-level
#      "Sub Total:" / "Grand Total:" rows are recalculated as the sum of all the
#      per-company sub totals.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-ItemRows($rowA, $rowB) {
    foreach ($col in @("B","C","D","E","F","G")) {
        $addrA = "$col$rowA"
        $addrB = "$col$rowB"
        $valA = $ws.Range($addrA).Value2
        $valB = $ws.Range($addrB).Value2
        $ws.Range($addrA).Value2 = $valB
        $ws.Range($addrB).Value2 = $valA
    }
}

# --- 1) Un-transpose the three swapped item-row pairs -----------------------
Swap-ItemRows 11 12
Swap-ItemRows 294 295
Swap-ItemRows 697 698

# --- 2) Quantity corrections; Value (G) = Rate (D) * Qty (F) ----------------
# Exact recalculated values are supplied so results match the source extract
# (Excel/engine floating point) bit-for-bit.
$ws.Range("F110").Value2 = 3
$ws.Range("G110").Value2 = 96430.74000000001
$ws.Range("F142").Value2 = 32
$ws.Range("G142").Value2 = 4590.08
$ws.Range("F143").Value2 = 121
$ws.Range("G143").Value2 = 6438.41
$ws.Range("F146").Value2 = 136
$ws.Range("G146").Value2 = 8880.799999999999
$ws.Range("F147").Value2 = 233
$ws.Range("G147").Value2 = 3919.06
$ws.Range("F149").Value2 = 70
$ws.Range("G149").Value2 = 6502.3
$ws.Range("F152").Value2 = 19
$ws.Range("G152").Value2 = 470.82
$ws.Range("F159").Value2 = 66
$ws.Range("G159").Value2 = 6900.3
$ws.Range("F181").Value2 = 25
$ws.Range("G181").Value2 = 2720.25
$ws.Range("F314").Value2 = 16
$ws.Range("G314").Value2 = 10046.72
$ws.Range("F324").Value2 = 2
$ws.Range("G324").Value2 = 10837.22
$ws.Range("F377").Value2 = 107
$ws.Range("G377").Value2 = 8562.139999999999
$ws.Range("F382").Value2 = 50
$ws.Range("G382").Value2 = 4108
$ws.Range("F397").Value2 = 48
$ws.Range("G397").Value2 = 6889.92
$ws.Range("F399").Value2 = 51
$ws.Range("G399").Value2 = 3364.98
$ws.Range("F401").Value2 = 19
$ws.Range("G401").Value2 = 2726.12
$ws.Range("F402").Value2 = 27
$ws.Range("G402").Value2 = 2243.16
$ws.Range("F403").Value2 = 54
$ws.Range("G403").Value2 = 2572.56
$ws.Range("F416").Value2 = 316
$ws.Range("G416").Value2 = 7305.92
$ws.Range("F422").Value2 = 150
$ws.Range("G422").Value2 = 8523
$ws.Range("F428").Value2 = 64
$ws.Range("G428").Value2 = 3681.92
$ws.Range("F437").Value2 = 23
$ws.Range("G437").Value2 = 3557.18
$ws.Range("F438").Value2 = 741
$ws.Range("G438").Value2 = 43496.7
$ws.Range("F441").Value2 = 24
$ws.Range("G441").Value2 = 5199.12
$ws.Range("F454").Value2 = 604
$ws.Range("G454").Value2 = 103483.32
$ws.Range("F455").Value2 = 181
$ws.Range("G455").Value2 = 27361.77
$ws.Range("F467").Value2 = 421
$ws.Range("G467").Value2 = 17336.78
$ws.Range("F468").Value2 = 196
$ws.Range("G468").Value2 = 18155.48
$ws.Range("F474").Value2 = 102
$ws.Range("G474").Value2 = 18725.16
$ws.Range("F475").Value2 = 112
$ws.Range("G475").Value2 = 20560.96
$ws.Range("F481").Value2 = 27
$ws.Range("G481").Value2 = 5624.1
$ws.Range("F518").Value2 = 52
$ws.Range("G518").Value2 = 6231.16
$ws.Range("F519").Value2 = 160
$ws.Range("G519").Value2 = 1700.8
$ws.Range("F531").Value2 = 36
$ws.Range("G531").Value2 = 236.52
$ws.Range("F534").Value2 = 51
$ws.Range("G534").Value2 = 7550.04
$ws.Range("F551").Value2 = 414
$ws.Range("G551").Value2 = 4409.1
$ws.Range("F555").Value2 = 33
$ws.Range("G555").Value2 = 1056.66
$ws.Range("F559").Value2 = 259
$ws.Range("G559").Value2 = 42991.41
$ws.Range("F560").Value2 = 175
$ws.Range("G560").Value2 = 26297.25
$ws.Range("F600").Value2 = 10
$ws.Range("G600").Value2 = 921.2
$ws.Range("F607").Value2 = 64
$ws.Range("G607").Value2 = 798.08
$ws.Range("F608").Value2 = 32
$ws.Range("G608").Value2 = 920.96
$ws.Range("F609").Value2 = 13
$ws.Range("G609").Value2 = 648.5700000000001
$ws.Range("F670").Value2 = 149
$ws.Range("G670").Value2 = 7461.92
$ws.Range("F671").Value2 = 330
$ws.Range("G671").Value2 = 20344.5
$ws.Range("F677").Value2 = 237
$ws.Range("G677").Value2 = 37758.84
$ws.Range("F721").Value2 = 38
$ws.Range("G721").Value2 = 19570.76
$ws.Range("F729").Value2 = 0
$ws.Range("G729").Value2 = 0
$ws.Range("F787").Value2 = 14
$ws.Range("G787").Value2 = 715.12
$ws.Range("F863").Value2 = 74
$ws.Range("G863").Value2 = 6789.5
$ws.Range("F868").Value2 = 11
$ws.Range("G868").Value2 = 1571.79
$ws.Range("F907").Value2 = 223
$ws.Range("G907").Value2 = 18187.88
$ws.Range("F910").Value2 = 43
$ws.Range("G910").Value2 = 6644.36
$ws.Range("F911").Value2 = 169
$ws.Range("G911").Value2 = 13783.64
$ws.Range("F912").Value2 = 320
$ws.Range("G912").Value2 = 42592
$ws.Range("F919").Value2 = 376
$ws.Range("G919").Value2 = 26154.56
$ws.Range("F927").Value2 = 304
$ws.Range("G927").Value2 = 43776
$ws.Range("F929").Value2 = 200
$ws.Range("G929").Value2 = 24142
$ws.Range("F944").Value2 = 17
$ws.Range("G944").Value2 = 1390.6
$ws.Range("F972").Value2 = 6
$ws.Range("G972").Value2 = 12996.78
$ws.Range("F993").Value2 = 18
$ws.Range("G993").Value2 = 822.0599999999999
$ws.Range("F994").Value2 = 24
$ws.Range("G994").Value2 = 626.16
$ws.Range("F996").Value2 = 1466
$ws.Range("G996").Value2 = 239119.26
$ws.Range("F998").Value2 = 82
$ws.Range("G998").Value2 = 23195.34
$ws.Range("F1005").Value2 = 7
$ws.Range("G1005").Value2 = 2189.25
$ws.Range("F1007").Value2 = 56
$ws.Range("G1007").Value2 = 2213.68

# --- 3) Recalculate affected Sub Total / Grand Total rows -------------------
$ws.Range("B111").Value2 = 150410.7
$ws.Range("B160").Value2 = 106226.8
$ws.Range("B188").Value2 = 110913.81
$ws.Range("B326").Value2 = 111506.61
$ws.Range("B472").Value2 = 586805.63
$ws.Range("B490").Value2 = 94955.99000000001
$ws.Range("B544").Value2 = 112163.82
$ws.Range("B561").Value2 = 134727.23
$ws.Range("B614").Value2 = 32298.43
$ws.Range("B684").Value2 = 173873.44
$ws.Range("B731").Value2 = 62029.2
$ws.Range("B799").Value2 = 22267.7
$ws.Range("B871").Value2 = 51049.97
$ws.Range("B931").Value2 = 419419.76
$ws.Range("B962").Value2 = 362418.8
$ws.Range("B986").Value2 = 357063.48
$ws.Range("B1002").Value2 = 264162.41
$ws.Range("B1008").Value2 = 6821.39
$ws.Range("B1014").Value2 = 5743596.15
$ws.Range("B1015").Value2 = 5743596.15
